$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the benchmark numbers in the "Throughput" table (C6:C15) and in
#    the "Bandwidth" table (C19:C28). These live on the worksheet and are the
#    source data the two charts read from.
# ---------------------------------------------------------------------------
$ws.Range("C6").Value  = 30.363322
$ws.Range("C7").Value  = 136.42182
$ws.Range("C8").Value  = 269.20534199999997
$ws.Range("C9").Value  = 444.94196299999999
$ws.Range("C10").Value = 1053.9363370000001
$ws.Range("C11").Value = 2145.7190639999999
$ws.Range("C12").Value = 4529.9939469999999
$ws.Range("C13").Value = 6451.8440959999998
$ws.Range("C14").Value = 6358.5518460000003
$ws.Range("C15").Value = 6314.6543220000003

$ws.Range("C19").Value = 0.54654000000000003
$ws.Range("C20").Value = 2.4555929999999999
$ws.Range("C21").Value = 4.8456960000000002
$ws.Range("C22").Value = 8.0089550000000003
$ws.Range("C23").Value = 18.970853999999999
$ws.Range("C24").Value = 38.622942999999999
$ws.Range("C25").Value = 81.539890999999997
$ws.Range("C26").Value = 116.133194
$ws.Range("C27").Value = 114.45393300000001
$ws.Range("C28").Value = 113.66377799999999

# ---------------------------------------------------------------------------
# 2) Update the chart titles: "Throughput M/s" -> "Throughput (M/s)" and
#    "Bandwidth GB/s" -> "Bandwidth (GB/s)".
# ---------------------------------------------------------------------------
$chart1 = $ws.ChartObjects().Item(1).Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Throughput (M/s)"

$chart2 = $ws.ChartObjects().Item(2).Chart
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Bandwidth (GB/s)"

# ---------------------------------------------------------------------------
# 3) Resize / reposition both chart objects on the sheet to their new
#    anchors (values derived from the target two-cell anchor positions,
#    converted to points: row height 15pt, column width 58.4375pt).
# ---------------------------------------------------------------------------
$co1 = $ws.ChartObjects().Item(1)
$co1.Top    = 7.12496062992126
$co1.Left   = 299.6874212598425
$co1.Width  = 697.5
$co1.Height = 257.62503937007875

$co2 = $ws.ChartObjects().Item(2)
$co2.Top    = 295.8748818897638
$co2.Left   = 304.1874212598425
$co2.Width  = 676.5625
$co2.Height = 229.8750393700787

# ---------------------------------------------------------------------------
# 4) Update the active cell selection on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("H20").Select()
